$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing benchmark values for the "32-32 / 17-20-3" block (rows 12-13).
# These look like numbers but must be stored as text (like the other timing
# cells already on this sheet, e.g. G8/H8/D9..H9). Writing the string straight
# into the destination cell makes Excel auto-convert it to a real number and
# changes the cell style (new NumberFormat xf). To avoid that, stage the text
# in a scratch cell formatted as Text, then copy/paste-special (values only)
# into the destination so the text type is kept but the destination's
# original style/number format is left untouched.
$cells = @("G12", "H12", "D13", "E13", "F13", "G13", "H13")
$values = @{
    "G12" = "4.51"
    "H12" = "4.50"
    "D13" = "5.27"
    "E13" = "3.18"
    "F13" = "3.20"
    "G13" = "8.43"
    "H13" = "8.42"
}

$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

foreach ($addr in $cells) {
    $scratch.Value = $values[$addr]
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}

$scratch.Clear()
$excel.CutCopyMode = $false

# Update the active selection to match the final cursor position
$ws.Range("H14").Select()
